# Adds the "Automated process / Manual process" legend diagram group to slide 1,
# matching the target OOXML diff (new p:grpSp "Group 5" with nested groups/shapes).
#
# NOTE on unit conversion: PowerPoint COM Left/Top/Width/Height are expressed in
# points (1 pt = 12700 EMU) and are stored as single-precision floats. To land on
# an exact target EMU value after the float round-trip (which truncates toward
# zero when converting back to EMU) we nudge the point value up by a small
# epsilon so the truncated result matches the desired EMU exactly.

function PtForEmu($targetEmu) {
    $basePt = $targetEmu / 12700.0
    for ($steps = 0; $steps -lt 2000; $steps++) {
        $candidate = $basePt + ($steps * 0.000001)
        # round-trip through a 32-bit float like the host engine does
        $f32 = [single]$candidate
        $emuOut = [int64]([double]$f32 * 12700.0)
        if ($emuOut -eq $targetEmu) {
            return $candidate
        }
    }
    return $basePt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Group 2: accent6 flowchart box + "Automated process" label
# ---------------------------------------------------------------------------

# Duplicate an existing accent6 "Flowchart: Alternate Process" shape so the
# p:style (lnRef/fillRef/effectRef accent6) comes across exactly.
$srcAccent6 = $s.Shapes.Item("Flowchart: Alternate Process 10")
$shape47 = $srcAccent6.Duplicate()
$shape47.Name = "Flowchart: Alternate Process 46"
$shape47.Left = PtForEmu 9787125
$shape47.Top = PtForEmu 1009169
$shape47.Width = PtForEmu 850416
$shape47.Height = PtForEmu 276999
$shape47.TextFrame.TextRange.Text = ""

$tb1 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb1.Name = "TextBox 1"
$tb1.Left = PtForEmu 10757770
$tb1.Top = PtForEmu 932224
$tb1.Width = PtForEmu 983411
$tb1.Height = PtForEmu 430887
$tb1.TextFrame.TextRange.Text = "Automated process"
$tb1.TextFrame.TextRange.Font.Size = 11
$tb1.TextFrame.WordWrap = -1
$tb1.TextFrame.AutoSize = 1

$range1 = $s.Shapes.Range(@($shape47.Name, $tb1.Name))
$group2 = $range1.Group()
$group2.Name = "Group 2"

# ---------------------------------------------------------------------------
# Group 4: accent1 flowchart box + "Manual process" label
# ---------------------------------------------------------------------------

$srcAccent1 = $s.Shapes.Item("Flowchart: Alternate Process 8")
$shape49 = $srcAccent1.Duplicate()
$shape49.Name = "Flowchart: Alternate Process 48"
$shape49.Left = PtForEmu 9800237
$shape49.Top = PtForEmu 1388141
$shape49.Width = PtForEmu 850416
$shape49.Height = PtForEmu 276999
$shape49.TextFrame.TextRange.Text = ""

$tb2 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb2.Name = "TextBox 49"
$tb2.Left = PtForEmu 10728387
$tb2.Top = PtForEmu 1311196
$tb2.Width = PtForEmu 983411
$tb2.Height = PtForEmu 430887
$tb2.TextFrame.TextRange.Text = "Manual process"
$tb2.TextFrame.TextRange.Font.Size = 11
$tb2.TextFrame.WordWrap = -1
$tb2.TextFrame.AutoSize = 1

$range2 = $s.Shapes.Range(@($shape49.Name, $tb2.Name))
$group4 = $range2.Group()
$group4.Name = "Group 4"

# ---------------------------------------------------------------------------
# Rectangle 50: bounding outline rectangle (duplicate of the big background
# rectangle so the line formatting/style match exactly)
# ---------------------------------------------------------------------------

$srcRect = $s.Shapes.Item("Rectangle 98")
$rect50 = $srcRect.Duplicate()
$rect50.Name = "Rectangle 50"
$rect50.Left = PtForEmu 9629910
$rect50.Top = PtForEmu 801153
$rect50.Width = PtForEmu 2196954
$rect50.Height = PtForEmu 1092470

# ---------------------------------------------------------------------------
# Group 5: outer group containing Group 2, Group 4 and Rectangle 50
# ---------------------------------------------------------------------------

$range3 = $s.Shapes.Range(@($group2.Name, $group4.Name, $rect50.Name))
$group5 = $range3.Group()
$group5.Name = "Group 5"

Write-Host "Diagram legend group added."
